# Update the "Förändrad" date column (C) for rows 2-8 from 2023-09-16 (45185)
# to 2023-10-05 (45204), keeping the existing date formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
